$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 3602631.33
$ws.Range("C9").Value = 568935.24
$ws.Range("D9").Value = 4171566.57
$ws.Range("E9").Value = 13.63840730941518
$ws.Range("F9").Value = 86.36159269058481
$ws.Range("G9").Value = -45.01524589588976
$ws.Range("H9").Value = -34.94145787994088
$ws.Range("I9").Value = 36353
$ws.Range("J9").Value = 1556
$ws.Range("K9").Value = 37909
$ws.Range("L9").Value = 26184
$ws.Range("M9").Value = 159.3173911549038
$ws.Range("N9").Value = 8.769227874966123
